$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1012
$ws1.Range("F4").Value = 13316
$ws1.Range("F6").Value = 1001
$ws1.Range("F7").Value = 5
$ws1.Range("F8").Value = 1589
$ws1.Range("F14").Value = 13303
$ws1.Range("F17").Value = 8879
$ws1.Range("F19").Value = 7957
$ws1.Range("F21").Value = 2
$ws1.Range("F23").Value = 427

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1012
$ws4.Range("F5").Value = 13316
$ws4.Range("F7").Value = 1001
$ws4.Range("F8").Value = 5
$ws4.Range("F9").Value = 1589
$ws4.Range("F15").Value = 13303
$ws4.Range("F18").Value = 8879
$ws4.Range("F20").Value = 7957
$ws4.Range("F22").Value = 2
$ws4.Range("F24").Value = 427

$wb.Save()
